$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update row 15: B15 1 -> 4, clear C15 (was 3) ---
$ws.Range("B15").Value = 4
$ws.Range("C15").ClearContents()

# --- Update row 16: B16 1.5 -> 2.75, C16 2.5 -> 1.25 ---
$ws.Range("B16").Value = 2.75
$ws.Range("C16").Value = 1.25

# --- Add new row 25 ---
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = 41462
$ws.Range("B25").Value = 2.5
$ws.Range("D25").Value = "SVN branch: gcc versus g++. Revision of Makefile, support of Linux and Windows, modularization"

# --- Update selection to new last cell ---
$ws.Range("A25").Select()
